$d = $word.ActiveDocument

# 1. Replace professor name
$d.Content.Find.Execute("Walter Priebe", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Joshua Sullivan", 2)

# 2. Add "or United States Space Force (USSF)" mention
$d.Content.Find.Execute("United States Air Force (USAF). The Air Force Reserve Officer Training Corps",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "United States Air Force (USAF) or United States Space Force (USSF). The Air Force Reserve Officer Training Corps",
                         2)

# 3. Add "the" before "Air Force Reserve Officer Training Corps" in AERO 1A/1B description
$d.Content.Find.Execute("encourages participation in Air Force Reserve Officer Training Corps.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "encourages participation in the Air Force Reserve Officer Training Corps.",
                         2)

# 4. Set page orientation to portrait explicitly on the section's page size
foreach ($sec in $d.Sections) {
    $sec.PageSetup.Orientation = 0
}
